$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C, rows 2 through 44 hold a "Förändrad" (changed) date.
# Bump each of these date values by one day (45180 -> 45181).
for ($r = 2; $r -le 44; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
